$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value corrections (imputed / corrected values) ---
# RM 8 (row 3): column C (D) was missing -> now has a value
$ws.Range("D3").Value = -14.2

# RM 9 (row 4): column D (E) had a value -> now missing
$ws.Range("E4").ClearContents()

# RM 14 (row 5): column C (D) had a value -> now missing
$ws.Range("D5").ClearContents()

# RM 42 (row 9): column D (E) was missing -> now has a value
$ws.Range("E9").Value = -6.8

# RM 52 a (row 10): column D (E) was missing -> now has a value
$ws.Range("E10").Value = -6.1

# RM 88 (row 13): column D (E) had a value -> now missing
$ws.Range("E13").ClearContents()

# RM 90 (row 14): column D (E) had a value -> now missing
$ws.Range("E14").ClearContents()

# RM 135 (row 21): column C (D) was missing -> now has a value
$ws.Range("D21").Value = -14.3

# RM 140 (row 23): column C (D) had a value -> now missing
$ws.Range("D23").ClearContents()

# --- Row removals ---
# Remove the "SC 92" row (originally row 28) first so the row-26 index
# used below still points at "RM 232".
$ws.Rows.Item(28).Delete()

# Remove the "RM 232" row (row 26)
$ws.Rows.Item(26).Delete()

# After the two deletions above, the row that used to be "SC 193" (row 34)
# is now row 32; its column C (D) value was missing -> now has a value.
$ws.Range("D32").Value = -14.7
